# Update cryptocurrency price/volume data and restore row order for Maker/VeChain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.303.59"
$ws.Range("E2").Value = "  +3.42%  "
$ws.Range("D3").Value = "'2.997.39"
$ws.Range("E3").Value = "  +3.55%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'562.90"
$ws.Range("E5").Value = "  +2.23%  "
$ws.Range("D6").Value = "'138.55"
$ws.Range("E6").Value = "  +12.36%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "'0.520"
$ws.Range("E8").Value = "  +4.75%  "
$ws.Range("D9").Value = "'2.986.62"
$ws.Range("E9").Value = "  +3.25%  "
$ws.Range("D10").Value = "'0.133"
$ws.Range("E10").Value = "  +7.52%  "
$ws.Range("D11").Value = "'4.89"
$ws.Range("E11").Value = "  +4.01%  "
$ws.Range("D12").Value = "'0.457"
$ws.Range("E12").Value = "  +4.03%  "
$ws.Range("D13").Value = "'0.0000232"
$ws.Range("E13").Value = "  +9.53%  "
$ws.Range("D14").Value = "'33.80"
$ws.Range("E14").Value = "  +3.98%  "
$ws.Range("E15").Value = "  +2.87%  "
$ws.Range("D16").Value = "'3.480.75"
$ws.Range("E16").Value = "  +3.17%  "
$ws.Range("D17").Value = "'7.01"
$ws.Range("E17").Value = "  +7.30%  "
$ws.Range("D18").Value = "'2.984.76"
$ws.Range("E18").Value = "  +3.03%  "
$ws.Range("D19").Value = "'59.071.69"
$ws.Range("E19").Value = "  +3.02%  "
$ws.Range("D20").Value = "'426.19"
$ws.Range("E20").Value = "  +5.36%  "
$ws.Range("D21").Value = "'13.55"
$ws.Range("E21").Value = "  +4.92%  "
$ws.Range("D22").Value = "'0.712"
$ws.Range("E22").Value = "  +6.06%  "
$ws.Range("D23").Value = "'7.17"
$ws.Range("E23").Value = "  +4.72%  "
$ws.Range("D24").Value = "'13.47"
$ws.Range("E24").Value = "  +5.62%  "
$ws.Range("D25").Value = "'80.67"
$ws.Range("E25").Value = "  +4.63%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").Value = "'0.998"
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("D28").Value = "'2.14"
$ws.Range("E28").Value = "  +10.24%  "
$ws.Range("D29").Value = "'2.54"
$ws.Range("E29").Value = "  +3.28%  "
$ws.Range("D30").Value = "'7.84"
$ws.Range("E30").Value = "  +9.16%  "
$ws.Range("D31").Value = "'25.66"
$ws.Range("E31").Value = "  +3.96%  "
$ws.Range("D32").Value = "'6.16"
$ws.Range("E32").Value = "  +2.21%  "
$ws.Range("D33").Value = "'0.0992"
$ws.Range("E33").Value = "  +0.50%  "
$ws.Range("E34").Value = "  +11.30%  "
$ws.Range("D35").Value = "'0.0₃0780"
$ws.Range("E35").Value = "  +25.77%  "
$ws.Range("D36").Value = "'5.76"
$ws.Range("E36").Value = "  +6.20%  "
$ws.Range("D37").Value = "'2.09"
$ws.Range("E37").Value = "  +4.74%  "
$ws.Range("D38").Value = "'49.00"
$ws.Range("E38").Value = "  +2.35%  "
$ws.Range("E39").Value = "  +3.78%  "
$ws.Range("D40").Value = "'2.83"
$ws.Range("E40").Value = "  +17.18%  "
$ws.Range("D41").Value = "'405.85"
$ws.Range("E41").Value = "  +12.82%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "'0.0351"
$ws.Range("E42").Value = "  +2.81%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "'2.743.45"
$ws.Range("E43").Value = "  +4.71%  "
$ws.Range("E44").Value = "  +1.55%  "
$ws.Range("D45").Value = "'0.247"
$ws.Range("E45").Value = "  +7.74%  "
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("D47").Value = "'125.13"
$ws.Range("E47").Value = "  +5.66%  "
$ws.Range("D48").Value = "'2.03"
$ws.Range("E48").Value = "  +4.41%  "
$ws.Range("E49").Value = "  +2.50%  "
$ws.Range("D50").Value = "'32.65"
$ws.Range("E50").Value = "  +20.48%  "
$ws.Range("D51").Value = "'23.51"
$ws.Range("E51").Value = "  +2.91%  "
